$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.048.43'
$ws.Range('E2').Value = '  -0.30%  '

$ws.Range('D3').Value = '2.304.35'
$ws.Range('E3').Value = '  -0.51%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '301.40'
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '98.71'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.94%  '

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.526'
$c.ClearFormats()
$ws.Range('E7').Value = '  +4.28%  '

$ws.Range('E8').Value = '  -0.04%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.524'
$c.ClearFormats()
$ws.Range('E9').Value = '  +1.33%  '

$ws.Range('E10').Value = '  -1.08%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0792'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.39%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.117'
$c.ClearFormats()
$ws.Range('E12').Value = '  -0.93%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '17.85'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.54%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.92'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.11%  '

$ws.Range('D15').Value = '2.663.67'
$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('D16').Value = '2.338.94'
$ws.Range('E16').Value = '  +0.85%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.789'
$c.ClearFormats()
$ws.Range('E17').Value = '  -2.49%  '

$ws.Range('D18').Value = '42.952.89'
$ws.Range('E18').Value = '  -0.29%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '13.48'
$c.ClearFormats()
$ws.Range('E19').Value = '  +6.71%  '

$ws.Range('E20').Value = '  +0.57%  '

$ws.Range('E21').Value = '  -0.73%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '68.31'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.65%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '239.40'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.13%  '

$ws.Range('E24').Value = '  -1.73%  '

$ws.Range('E25').Value = '  -0.02%  '

$ws.Range('E26').Value = '  -1.06%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '24.83'
$c.ClearFormats()
$ws.Range('E27').Value = '  +0.25%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '167.40'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.89%  '

$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '9.13'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.89%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.04'
$c.ClearFormats()
$ws.Range('E30').Value = '  -13.38%  '

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '33.36'
$c.ClearFormats()
$ws.Range('E31').Value = '  -4.20%  '

$ws.Range('E32').Value = '  +4.19%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E33').Value = '  -0.04%  '

$ws.Range('E34').Value = '  +2.01%  '

$ws.Range('E35').Value = '  +4.41%  '

$ws.Range('E36').Value = '  -0.51%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.0690'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.68%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range('E38').Value = '  -0.95%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '1.80'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.05%  '

$ws.Range('E40').Value = '  +1.98%  '

$ws.Range('E41').Value = '  -3.28%  '

$ws.Range('D42').Value = '2.005.02'
$ws.Range('E42').Value = '  +0.80%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.0288'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.51%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '10.10'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.38%  '

$ws.Range('E45').Value = '  -5.62%  '

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '17.37'
$c.ClearFormats()
$ws.Range('E46').Value = '  -1.23%  '

$ws.Range('E47').Value = '  -2.95%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '54.61'
$c.ClearFormats()
$ws.Range('E48').Value = '  -2.89%  '

$ws.Range('D49').Value = '2.529.42'
$ws.Range('E49').Value = '  -0.68%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.54'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.28%  '

$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '73.49'
$c.ClearFormats()
$ws.Range('E51').Value = '  +4.66%  '
